# Auto-generated edit script applying numeric corrections to the
# Zodiark_Profits workbook (per-sheet currentAveragePrice / LevePrice /
# LeveProfit recalculations from the scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 3906.6538
$ws.Range("I2").Value = 576
$ws.Range("J2").Value = 5988.3125
$ws.Range("K2").Value = 576
$ws.Range("L2").Value = 5988.3125
$ws.Range("M2").Value = -463
$ws.Range("N2").Value = -6214.3125
$ws.Range("H9").Value = 120.333336
$ws.Range("I9").Value = 127
$ws.Range("K9").Value = 127
$ws.Range("M9").Value = 42
$ws.Range("H62").Value = 8927.799999999999
$ws.Range("I62").Value = 8927.799999999999
$ws.Range("K62").Value = 8927.799999999999
$ws.Range("M62").Value = -8303.799999999999
$ws.Range("H65").Value = 8927.799999999999
$ws.Range("I65").Value = 8927.799999999999
$ws.Range("K65").Value = 44639
$ws.Range("M65").Value = -41519
$ws.Range("H86").Value = 1807.3636
$ws.Range("I86").Value = 1948.5
$ws.Range("K86").Value = 1948.5
$ws.Range("M86").Value = -825.5
$ws.Range("H89").Value = 1807.3636
$ws.Range("I89").Value = 1948.5
$ws.Range("K89").Value = 9742.5
$ws.Range("M89").Value = -4126.5
$ws.Range("H111").Value = 1268
$ws.Range("I111").Value = 937
$ws.Range("K111").Value = 2811
$ws.Range("M111").Value = 256
$ws.Range("H132").Value = 2197.5386
$ws.Range("I132").Value = 2213.6667
$ws.Range("K132").Value = 6641.000100000001
$ws.Range("M132").Value = -4111.000100000001
$ws.Range("H138").Value = 2144.4102
$ws.Range("J138").Value = 2868.7727
$ws.Range("L138").Value = 8606.3181
$ws.Range("N138").Value = -18886.3181

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 3952
$ws.Range("J46").Value = 3952
$ws.Range("L46").Value = 3952
$ws.Range("N46").Value = -4590
$ws.Range("H122").Value = 7240.9
$ws.Range("I122").Value = 8265.643
$ws.Range("K122").Value = 24796.929
$ws.Range("M122").Value = -22346.929
$ws.Range("H133").Value = 83318.836
$ws.Range("J133").Value = 83318.836
$ws.Range("L133").Value = 83318.836
$ws.Range("N133").Value = -88378.836

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 849
$ws.Range("I11").Value = 849
$ws.Range("K11").Value = 849
$ws.Range("M11").Value = -709
$ws.Range("H24").Value = 3138.6667
$ws.Range("J24").Value = 2250
$ws.Range("L24").Value = 2250
$ws.Range("N24").Value = -2720
$ws.Range("H31").Value = 7905
$ws.Range("I31").Value = 7905
$ws.Range("K31").Value = 7905
$ws.Range("M31").Value = -7653

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 25.333334
$ws.Range("I7").Value = 27.4
$ws.Range("J7").Value = 22.75
$ws.Range("K7").Value = 27.4
$ws.Range("L7").Value = 22.75
$ws.Range("M7").Value = 85.59999999999999
$ws.Range("N7").Value = -248.75
$ws.Range("H31").Value = 1473.5
$ws.Range("I31").Value = 1083.1777
$ws.Range("J31").Value = 2506.7058
$ws.Range("K31").Value = 1083.1777
$ws.Range("L31").Value = 2506.7058
$ws.Range("M31").Value = -788.1777
$ws.Range("N31").Value = -3096.7058
$ws.Range("H34").Value = 1473.5
$ws.Range("I34").Value = 1083.1777
$ws.Range("J34").Value = 2506.7058
$ws.Range("K34").Value = 1083.1777
$ws.Range("L34").Value = 2506.7058
$ws.Range("M34").Value = -881.1777
$ws.Range("N34").Value = -2910.7058
$ws.Range("H86").Value = 166672770
$ws.Range("I86").Value = 166672770
$ws.Range("K86").Value = 166672770
$ws.Range("M86").Value = -166671647
$ws.Range("H88").Value = 15742
$ws.Range("J88").Value = 15742
$ws.Range("L88").Value = 15742
$ws.Range("N88").Value = -16554
$ws.Range("H89").Value = 166672770
$ws.Range("I89").Value = 166672770
$ws.Range("K89").Value = 833363850
$ws.Range("M89").Value = -833358234
$ws.Range("H91").Value = 15742
$ws.Range("J91").Value = 15742
$ws.Range("L91").Value = 15742
$ws.Range("N91").Value = -18550
$ws.Range("H107").Value = 896.5
$ws.Range("J107").Value = 1047.037
$ws.Range("L107").Value = 1047.037
$ws.Range("N107").Value = -4887.037
$ws.Range("H122").Value = 999
$ws.Range("I122").Value = 999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2997
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -547
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1999.8
$ws.Range("I132").Value = 1999.8
$ws.Range("K132").Value = 5999.4
$ws.Range("M132").Value = -3469.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 11399.5
$ws.Range("J104").Value = 11399.5
$ws.Range("L104").Value = 34198.5
$ws.Range("N104").Value = -39440.5
$ws.Range("H116").Value = 2992.3333
$ws.Range("I116").Value = 2994
$ws.Range("J116").Value = 2989
$ws.Range("K116").Value = 8982
$ws.Range("L116").Value = 8967
$ws.Range("M116").Value = -5540
$ws.Range("N116").Value = -15851
$ws.Range("H132").Value = 1886.0714
$ws.Range("I132").Value = 1868
$ws.Range("J132").Value = 1896.1111
$ws.Range("K132").Value = 16812
$ws.Range("L132").Value = 17064.9999
$ws.Range("M132").Value = -14282
$ws.Range("N132").Value = -22124.9999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 25507.424
$ws.Range("J70").Value = 6320.3335
$ws.Range("L70").Value = 6320.3335
$ws.Range("N70").Value = -6860.3335
$ws.Range("H73").Value = 25507.424
$ws.Range("J73").Value = 6320.3335
$ws.Range("L73").Value = 6320.3335
$ws.Range("N73").Value = -8192.333500000001
$ws.Range("H80").Value = 2965.516
$ws.Range("I80").Value = 2786.8572
$ws.Range("J80").Value = 3340.7
$ws.Range("K80").Value = 2786.8572
$ws.Range("L80").Value = 3340.7
$ws.Range("M80").Value = -1788.8572
$ws.Range("N80").Value = -5336.7
$ws.Range("H83").Value = 2965.516
$ws.Range("I83").Value = 2786.8572
$ws.Range("J83").Value = 3340.7
$ws.Range("K83").Value = 13934.286
$ws.Range("L83").Value = 16703.5
$ws.Range("M83").Value = -8942.286
$ws.Range("N83").Value = -26687.5
$ws.Range("H102").Value = 2618.842
$ws.Range("I102").Value = 2581.7856
$ws.Range("K102").Value = 2581.7856
$ws.Range("M102").Value = -959.7856000000002
$ws.Range("H122").Value = 1990.3793
$ws.Range("I122").Value = 1379.0526
$ws.Range("J122").Value = 3151.9
$ws.Range("K122").Value = 4137.1578
$ws.Range("L122").Value = 9455.700000000001
$ws.Range("M122").Value = -1687.1578
$ws.Range("N122").Value = -14355.7
$ws.Range("H126").Value = 5136530.5
$ws.Range("I126").Value = 4050.8572
$ws.Range("K126").Value = 12152.5716
$ws.Range("M126").Value = -9682.571599999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4564.7095
$ws.Range("I7").Value = 4455.125
$ws.Range("K7").Value = 4455.125
$ws.Range("M7").Value = -4343.125
$ws.Range("H22").Value = 5099.375
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 5685.143
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 5685.143
$ws.Range("M22").Value = -704
$ws.Range("N22").Value = -6275.143
$ws.Range("H27").Value = 5099.375
$ws.Range("I27").Value = 999
$ws.Range("J27").Value = 5685.143
$ws.Range("K27").Value = 999
$ws.Range("L27").Value = 5685.143
$ws.Range("M27").Value = -892
$ws.Range("N27").Value = -5899.143
$ws.Range("H40").Value = 4245.85
$ws.Range("I40").Value = 3557.4707
$ws.Range("J40").Value = 8146.6665
$ws.Range("K40").Value = 3557.4707
$ws.Range("L40").Value = 8146.6665
$ws.Range("M40").Value = -3421.4707
$ws.Range("N40").Value = -8418.666499999999
$ws.Range("H122").Value = 5533.875
$ws.Range("I122").Value = 5429.1763
$ws.Range("K122").Value = 16287.5289
$ws.Range("M122").Value = -13837.5289
$ws.Range("H126").Value = 4564.7095
$ws.Range("I126").Value = 4455.125
$ws.Range("K126").Value = 13365.375
$ws.Range("M126").Value = -10895.375
$ws.Range("H132").Value = 4917.778
$ws.Range("I132").Value = 4742.923
$ws.Range("J132").Value = 5372.4
$ws.Range("K132").Value = 14228.769
$ws.Range("L132").Value = 16117.2
$ws.Range("M132").Value = -11698.769
$ws.Range("N132").Value = -21177.2
$ws.Range("H136").Value = 1819.12
$ws.Range("I136").Value = 1478.0416
$ws.Range("K136").Value = 4434.1248
$ws.Range("M136").Value = -1884.1248

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 536
$ws.Range("I100").Value = 317.77777
$ws.Range("K100").Value = 635.55554
$ws.Range("M100").Value = -94.55553999999995
$ws.Range("H113").Value = 2640.1667
$ws.Range("I113").Value = 2368.2
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 7104.599999999999
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = -4934.599999999999
$ws.Range("N113").Value = -16340
$ws.Range("H126").Value = 142859740
$ws.Range("J126").Value = 500001000
$ws.Range("L126").Value = 1500003000
$ws.Range("N126").Value = -1500007940
$ws.Range("H132").Value = 2153
$ws.Range("I132").Value = 1974.2858
$ws.Range("J132").Value = 2570
$ws.Range("K132").Value = 5922.857400000001
$ws.Range("L132").Value = 7710
$ws.Range("M132").Value = -3392.857400000001
$ws.Range("N132").Value = -12770
